$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Clear the tries/runtime/memory/url/note data for row 9
#    ("Remove Duplicates from Sorted Array") - only day/problemName remain.
$ws1.Range("C9:H9").ClearContents()

# 2. Add new day-15 entry: "Climbing Stairs" with a "fibonacci " note.
$ws1.Range("A17").Value = 15
$ws1.Range("B17").Value = "Climbing Stairs"
$ws1.Range("I17").Value = "fibonacci "

# 3. Add the remaining day numbers (16 through 33) in column A only.
for ($r = 18; $r -le 35; $r++) {
    $ws1.Range("A$r").Value = $r - 2
}

# 4. Widen column B (problemName) to fit the new long entry; drop best-fit.
$ws1.Columns.Item(2).ColumnWidth = 37.14

# 5. Update the saved selections: Sheet2's selection moves to C30 (set
#    first, without leaving Sheet2 as the active tab), then Sheet1's
#    selection moves to H24 and Sheet1 stays the active sheet.
$ws2.Range("C30").Select()
$ws1.Activate()
$ws1.Range("H24").Select()
